$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1 ("About")
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("About")

# Update the "As of EPS ..." explanatory paragraph (rows 11-13) and
# add a new 4th line of text (row 14), then drop the now-obsolete
# paragraph that used to live in rows 15-18.
$ws1.Range("A11").Value = "As of EPS 1.5.0, this lever supports the three energy carriers (electricity,"
$ws1.Range("A12").Value = "district heat, and hydrogen), which tend to be produced and consumed locally."
$ws1.Range("A13").Value = "It does not affect other fuel types, whose prices are often determined or influenced"
$ws1.Range("A14").Value = "by global markets, so domestic producers' costs are less relevant."
$ws1.Range("A15:A18").ClearContents()

# ---------------------------------------------------------------
# Sheet 2 ("BAEPAbCiPC")
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("BAEPAbCiPC")

# Header cell: shorten label and drop the italic styling
$ws2.Range("A1").Value = "Boolean"
$ws2.Range("A1").Style = "Normal"

# "nuclear" row loses the "(NOT USED)" suffix and its grey shading
$ws2.Range("A5").Value = "nuclear"
$ws2.Range("A5").Style = "Normal"
$ws2.Range("B5").Style = "Normal"

# "municipal solid waste" row loses the "(NOT USED)" suffix and its grey shading
$ws2.Range("A21").Value = "municipal solid waste"
$ws2.Range("A21").Style = "Normal"
$ws2.Range("B21").Style = "Normal"

# The three energy-carrier rows (electricity, heat, hydrogen) swap their
# italic styling for a green highlight fill (RGB 146, 208, 80 = #92D050)
$greenFill = 5296274  # OLE (BGR) value of RGB(146, 208, 80)
foreach ($row in 2, 15, 22) {
    $ws2.Range("A$row").Style = "Normal"
    $ws2.Range("B$row").Style = "Normal"
    $ws2.Range("A$row").Interior.Color = $greenFill
    $ws2.Range("B$row").Interior.Color = $greenFill
}
